# Advection tests with smaller dt
# Add two new test-case rows (19 and 20) below the existing "swim to deep,
# mult by lmask 1st" row, mirroring the formatting of the rows immediately
# above them, then move the active selection to A21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: "swim to const rand" -----------------------------------
# Copy row 18's formatting (label style + value styles) down to row 19.
$ws.Range("A18:I18").Copy()
$ws.Range("A19:I19").PasteSpecial(-4122)
$ws.Rows.Item(19).RowHeight = 32

$ws.Range("A19").Value = "New 200m avg vel, dt = 1 hr, j = 2, swim to const rand, mult by lmask 1st"
$ws.Range("C19").Value = [double]"2.6601999999999998E+37"
$ws.Range("E19").Value = [double]"3.1471000000000002E+46"
$ws.Range("I19").Value = [double]"1.5424999999999999E+42"
$ws.Range("C19").NumberFormat = "0.00E+00"
$ws.Range("I19").NumberFormat = "0.00E+00"

# --- Row 20: "swim to changing rand" ---------------------------------
# Copy row 17's formatting (plain General-format values) down to row 20.
$ws.Range("A17:I17").Copy()
$ws.Range("A20:I20").PasteSpecial(-4122)
$ws.Rows.Item(20).RowHeight = 32

$ws.Range("A20").Value = "New 200m avg vel, dt = 1 hr, j = 2, swim to changing rand, mult by lmask 1st"
$ws.Range("C20").Value = [double]"4.7473000000000001"
$ws.Range("E20").Value = [double]"161.745"
$ws.Range("I20").Value = [double]"177.05170000000001"

$excel.CutCopyMode = $false

# --- Move the active selection, matching the author's final cursor pos
$ws.Range("A21").Select() | Out-Null
